$d = $word.ActiveDocument

# --- Paragraph 1: title ---
$p1 = $d.Paragraphs(1)
$full1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">On</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">Pilgrimage</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">-</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">July</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">/</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">August</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">1977</w:t></w:r>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full1.InsertXML($xml1)
$p1.Style = "Title"

# --- Paragraph 2: authors ---
$p2 = $d.Paragraphs(2)
$full2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
  '<w:r><w:t xml:space="preserve">Dorothy</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">Day</w:t></w:r>' +
  '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$full2.InsertXML($xml2)
$p2.Style = "Authors"
